$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Determine the last used row in the sheet (data runs from row 2 to 481).
$lastRow = $ws.UsedRange.Rows.Count + $ws.UsedRange.Row - 1
if ($lastRow -lt 2) { $lastRow = 481 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45205) {
        $cell.Value2 = 45206
    }
}
